# Append a new data row (row 3) to the HZNPRandom sheet, mirroring the
# existing "Random" sample row (row 2) that was generated by the trader.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Duplicate row 2's formatting (keeps the date style on column A) into row 3,
# then overwrite the values with the new sample's data.
$ws.Range("A2:N2").Copy($ws.Range("A3:N3"))

$ws.Range("A3").Value = 42605.88585648148
$ws.Range("B3").Value = 4
$ws.Range("C3").Value = 0
$ws.Range("D3").Value = 0
$ws.Range("E3").Value = 0
$ws.Range("F3").Value = 0
$ws.Range("G3").Value = 0
$ws.Range("H3").Value = 0
$ws.Range("I3").Value = 0
$ws.Range("J3").Value = 0
$ws.Range("K3").Value = 0
$ws.Range("L3").Value = 0
$ws.Range("M3").Value = 0
$ws.Range("N3").Value = "Random"
